# Apply the "fix download docx file" edit to Statistica.docx
#
# Changes:
#   1) Heading "2 взвод"                       -> "31 курс"
#   2) Period text (paragraph)  "11.05.2023 – 09.06.2023" -> "12.05.2023 – 10.06.2023"
#   3) Period text (table cell) "11.05.2023 – 09.06.2023" -> "12.05.2023 – 10.06.2023"
#   4) Table cell "Всего"      value "2" -> "3"
#   5) Table cell "Взыскания"  value "0" -> "1"

$d = $word.ActiveDocument

# 1) Replace the plateau/group heading "2 взвод" with "31 курс"
$d.Content.Find.Execute(
    "2 взвод", $true, $false, $false, $false, $false,
    $true, 1, $false, "31 курс", 2
)

# 2) & 3) Replace every occurrence of the reporting-period date range.
#    Both the standalone paragraph and the table cell use the identical
#    text, so a global Find/Replace handles both occurrences.
$oldRange = "11.05.2023 " + [char]0x2013 + " 09.06.2023"
$newRange = "12.05.2023 " + [char]0x2013 + " 10.06.2023"
$d.Content.Find.Execute(
    $oldRange, $true, $false, $false, $false, $false,
    $true, 1, $false, $newRange, 2
)

# 4) & 5) Update the numeric totals in the data row of the statistics
#    table. The single table on the page has a header row followed by
#    one data row: Период | Всего | Поощрения | Взыскания | Снятия взыскания.
$table = $d.Tables.Item(1)

# "Всего" column (2nd column), data row (2nd row): "2" -> "3"
$table.Cell(2, 2).Range.Text = "3"

# "Взыскания" column (4th column), data row (2nd row): "0" -> "1"
$table.Cell(2, 4).Range.Text = "1"

$d.Save()
